# Regenerate orders with updated distance/sizes.
# The experiment's Distance (D51/D64/D80) and Size (S30) condition codes
# were renumbered (D51->D55, D64->D69, D80->D86, S30->S31). Every string
# in the workbook that encodes these tokens (condition labels, filenames,
# and the Distance/Size lookup lists) needs the same substitution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters only in that none of the new tokens collide with any
# other existing token, so a straightforward sequential replace is safe.
$null = $ws.Cells.Replace("D51", "D55")
$null = $ws.Cells.Replace("D64", "D69")
$null = $ws.Cells.Replace("D80", "D86")
$null = $ws.Cells.Replace("S30", "S31")
